$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = [double]"9.611007164012482E-16"
$ws.Range("E2").Value = [double]"9.611007164012482E-16"

$ws.Range("D3").Value = [double]"1.668103725087027E-09"
$ws.Range("E3").Value = [double]"1.668103725087027E-09"

$ws.Range("D4").Value = 0.9019478918766364
$ws.Range("E4").Value = 0.9019478918766364

$ws.Range("D5").Value = [double]"8.927797162911799E-30"
$ws.Range("E5").Value = [double]"8.927797162911799E-30"

$ws.Range("D6").Value = [double]"3.979863960769505E-32"
$ws.Range("E6").Value = [double]"3.979863960769505E-32"

$ws.Range("D7").Value = 0.9999999999994784
$ws.Range("E7").Value = [double]"5.215827769688985E-13"

$ws.Range("D8").Value = 0.9999999999999876
$ws.Range("E8").Value = [double]"1.243449787580175E-14"

$ws.Range("D10").Value = 0.4974537331742029
$ws.Range("E10").Value = 0.5025462668257972

$ws.Range("D11").Value = 0.9173177702926415
$ws.Range("E11").Value = 0.0826822297073585
$ws.Range("F11").Value = 87.40620422363281
